$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Miss Dina Nasr, Administrator") {
        $cell.Value = "Administrator, Miss Dina Nasr"
    }
}
